# Update computed balance (center-of-gravity) results on the
# "GLOBAL RESULTS" and "LANDING GEARS" sheets following a refresh of the
# performance analysis (assigned polar curve used for analyses).

$wb = $excel.ActiveWorkbook

# --- GLOBAL RESULTS sheet ---------------------------------------------
$wsGlobal = $wb.Worksheets.Item("GLOBAL RESULTS")

$wsGlobal.Range("C2").Value = 5.236520021021633
$wsGlobal.Range("C3").Value = 12.240997717128838
$wsGlobal.Range("C4").Value = 0.6576665053638682

$wsGlobal.Range("C6").Value = 5.041226939529571
$wsGlobal.Range("C7").Value = 11.791351849019662
$wsGlobal.Range("C8").Value = 0.721704678058392

$wsGlobal.Range("C10").Value = 5.041226939529571
$wsGlobal.Range("C11").Value = 11.791351849019662
$wsGlobal.Range("C12").Value = 0.721704678058392

$wsGlobal.Range("C14").Value = 5.081355335965858
$wsGlobal.Range("C15").Value = 11.883744101287732
$wsGlobal.Range("C16").Value = 0.46534746772689106

$wsGlobal.Range("C18").Value = 4.974961100839947
$wsGlobal.Range("C19").Value = 11.638780337853301
$wsGlobal.Range("C20").Value = 0.6736251083832776

# --- LANDING GEARS sheet ------------------------------------------------
$wsLandingGears = $wb.Worksheets.Item("LANDING GEARS")

$wsLandingGears.Range("C2").Value = 12.299024241711926
